$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "353.79", "23.80") keep their exact text representation instead
# of being auto-parsed into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '52.010.87'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '2.967.02'
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '353.79'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '112.09'
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").Value = '0.567'
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("E11").Value = '  +5.18%  '
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").Value = '19.97'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").Value = '8.06'
$ws.Range("E14").Value = '  +2.59%  '
$ws.Range("D15").Value = '3.431.95'
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").Value = '2.967.85'
$ws.Range("E16").Value = '  +2.13%  '
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").Value = '52.111.22'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '14.57'
$ws.Range("E20").Value = '  +6.17%  '
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("D22").Value = '0.0₃0992'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = '271.28'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").Value = '0.182'
$ws.Range("E26").Value = '  +10.19%  '
$ws.Range("D27").Value = '27.51'
$ws.Range("E27").Value = '  +3.34%  '
$ws.Range("D28").Value = '7.57'
$ws.Range("E28").Value = '  +19.74%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = '0.109'
$ws.Range("E30").Value = '  +21.30%  '
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("D32").Value = '37.83'
$ws.Range("E32").Value = '  -2.84%  '
$ws.Range("D33").Value = '6.26'
$ws.Range("E33").Value = '  +9.96%  '
$ws.Range("D34").Value = '53.19'
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("D35").Value = '0.0453'
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").Value = '1.97'
$ws.Range("E36").Value = '  -13.49%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '3.43'
$ws.Range("E38").Value = '  +2.93%  '
$ws.Range("D39").Value = '19.05'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  +1.41%  '
$ws.Range("D41").Value = '2.71'
$ws.Range("E41").Value = '  +4.20%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.119'
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '23.80'
$ws.Range("E43").Value = '  +4.65%  '
$ws.Range("E44").Value = '  -2.26%  '
$ws.Range("D45").Value = '3.58'
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("D47").Value = '2.182.08'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").Value = '113.79'
$ws.Range("E48").Value = '  -7.18%  '
$ws.Range("D49").Value = '0.243'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").Value = '0.0342'
$ws.Range("E50").Value = '  +6.46%  '
$ws.Range("D51").Value = '0.941'
$ws.Range("E51").Value = '  -2.49%  '
